$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their text representation (values like
# "0.9998" or "5.450" would otherwise be auto-converted to numbers, losing
# trailing zeros / exact formatting), by pre-formatting the column as Text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.815.74'
$ws.Range("E2").Value = '  -1.33%  '
$ws.Range("D3").Value = '1.892.62'
$ws.Range("E3").Value = '  -1.11%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '0.7769'
$ws.Range("E5").Value = '  -5.59%  '
$ws.Range("D6").Value = '243.87'
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("D8").Value = '0.3125'
$ws.Range("E8").Value = '  -4.17%  '
$ws.Range("D9").Value = '25.28'
$ws.Range("E9").Value = '  -7.21%  '
$ws.Range("D10").Value = '0.07174'
$ws.Range("E10").Value = '  +0.77%  '
$ws.Range("D11").Value = '0.08059'
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("D12").Value = '0.7653'
$ws.Range("E12").Value = '  -1.89%  '
$ws.Range("D13").Value = '5.450'
$ws.Range("E13").Value = '  +1.67%  '
$ws.Range("D14").Value = '1.861.12'
$ws.Range("E14").Value = '  -4.38%  '
$ws.Range("D15").Value = '92.26'
$ws.Range("E15").Value = '  -2.40%  '
$ws.Range("D16").Value = '6.166'
$ws.Range("D17").Value = '29.790.66'
$ws.Range("E17").Value = '  -1.51%  '
$ws.Range("D18").Value = '13.95'
$ws.Range("E18").Value = '  -2.81%  '
$ws.Range("D19").Value = '243.31'
$ws.Range("E19").Value = '  -2.32%  '
$ws.Range("D20").Value = '0.000007759'
$ws.Range("E20").Value = '  -1.08%  '
$ws.Range("E21").Value = '  -0.12%  '
$ws.Range("D22").Value = '8.099'
$ws.Range("E22").Value = '  +6.16%  '
$ws.Range("D23").Value = '2.106.63'
$ws.Range("E23").Value = '  -4.30%  '
$ws.Range("D24").Value = '1.000'
$ws.Range("D25").Value = '0.1604'
$ws.Range("E25").Value = '  -4.75%  '
$ws.Range("D26").Value = '9.398'
$ws.Range("E26").Value = '  -0.63%  '
$ws.Range("D27").Value = '161.50'
$ws.Range("E27").Value = '  -4.06%  '
$ws.Range("D28").Value = '18.73'
$ws.Range("E28").Value = '  -2.06%  '
$ws.Range("D29").Value = '2.047'
$ws.Range("E29").Value = '  -3.49%  '
$ws.Range("D30").Value = '1.431'
$ws.Range("E30").Value = '  +4.61%  '
$ws.Range("E31").Value = '  +0.94%  '
$ws.Range("D32").Value = '4.471'
$ws.Range("E32").Value = '  +2.93%  '
$ws.Range("E33").Value = '  -0.87%  '
$ws.Range("D34").Value = '0.05530'
$ws.Range("E34").Value = '  -2.63%  '
$ws.Range("D35").Value = '1.264'
$ws.Range("E35").Value = '  -1.89%  '
$ws.Range("D36").Value = '0.7451'
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("D37").Value = '0.9958'
$ws.Range("E37").Value = '  -0.42%  '
$ws.Range("D38").Value = '2.617'
$ws.Range("E38").Value = '  -3.56%  '
$ws.Range("D39").Value = '0.01913'
$ws.Range("E39").Value = '  -1.43%  '
$ws.Range("D40").Value = '2.778'
$ws.Range("E40").Value = '  -1.22%  '
$ws.Range("D41").Value = '1.139.87'
$ws.Range("E41").Value = '  +10.01%  '
$ws.Range("D42").Value = '73.59'
$ws.Range("E42").Value = '  -0.56%  '
$ws.Range("D43").Value = '0.4419'
$ws.Range("E43").Value = '  -1.45%  '
$ws.Range("D44").Value = '5.853'
$ws.Range("E44").Value = '  -2.16%  '
$ws.Range("D45").Value = '0.8531'
$ws.Range("E45").Value = '  +0.43%  '
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").Value = '103.74'
$ws.Range("E47").Value = '  +0.76%  '
$ws.Range("D48").Value = '1.887'
$ws.Range("E48").Value = '  -2.41%  '
$ws.Range("D49").Value = '9.914'
$ws.Range("E49").Value = '  -0.45%  '
$ws.Range("D50").Value = '7.441'
$ws.Range("E50").Value = '  -2.32%  '
$ws.Range("D51").Value = '3.013'
$ws.Range("E51").Value = '  +10.34%  '
